# Insert Marielle's IP address into the access spreadsheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marielle Miziara's row (row 11) was missing the IP column (D); fill it in.
$ws.Range("D11").Value = "192.168.0.120"

# Move the active selection to reflect where the user's cursor ended up
# after making the edit (one row below the last data row).
$ws.Range("D12").Select() | Out-Null
